$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 11 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 11
}

# New Actual Production (MW) values for rows 2-41 (column B)
$newValues = @{
    2  = 284
    3  = 255
    4  = 229
    5  = 228
    6  = 230
    7  = 221
    8  = 220
    9  = 222
    10 = 245
    11 = 255
    12 = 266
    13 = 272
    14 = 277
    15 = 255
    16 = 250
    17 = 244
    18 = 231
    19 = 227
    20 = 219
    21 = 218
    22 = 219
    23 = 229
    24 = 241
    25 = 232
    26 = 222
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
